$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Neutrophils"-originated row becomes MuSCs/Calca/Calcr/ECs with
# refreshed TPM-derived metrics.
$ws.Range("A2").Value = "MuSCs"
$ws.Range("B2").Value = "Calca"
$ws.Range("C2").Value = "Calcr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.347953
$ws.Range("H2").Value = 1.043859
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.09584066666666667
$ws.Range("N2").Value = 0.287522
$ws.Range("O2").Value = 0.3176649051884244
$ws.Range("P2").Value = 0.3176649051884244
$ws.Range("Q2").Value = 0.03334804748866667
$ws.Range("R2").Value = 0.300132427398
$ws.Range("S2").Value = 0.3176649051884244
$ws.Range("T2").Value = 0.3176649051884244

# Row 3: MuSCs/Calca/Calcr/MuSCs with refreshed TPM-derived metrics.
$ws.Range("A3").Value = "MuSCs"
$ws.Range("B3").Value = "Calca"
$ws.Range("C3").Value = "Calcr"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.347953
$ws.Range("H3").Value = 1.043859
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.205863
$ws.Range("N3").Value = 0.6175889999999999
$ws.Range("O3").Value = 0.6823350948115756
$ws.Range("P3").Value = 0.6823350948115755
$ws.Range("Q3").Value = 0.07163064843899999
$ws.Range("R3").Value = 0.644675835951
$ws.Range("S3").Value = 0.6823350948115756
$ws.Range("T3").Value = 0.6823350948115755
